$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers (human readable labels)
$ws.Range("A1").Value = "Territorio:"
$ws.Range("B1").Value = "Sector descripción"
$ws.Range("C1").Value = "Sector código"
$ws.Range("D1").Value = "Codcom"
$ws.Range("E1").Value = "Comarca nombre"
$ws.Range("F1").Value = "Dirección provincial nombre"
$ws.Range("G1").Value = "Mes y año"
$ws.Range("H1").Value = "Cuentas cotización con trabajadores"
$ws.Range("I1").Value = "Dirección provincial (código)"

# Row 2 - measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:territorio"
$ws.Range("B2").Value = "iaest-measure:sector-descripcion"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:codcom"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("G2").Value = "iaest-measure:mes-y-ano"
$ws.Range("H2").Value = "iaest-measure:cuentas-cotizacion-con-trabajadores"
$ws.Range("I2").Value = "null"

# Row 3 - medida/dim/null markers
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "null"

# Row 4 - xsd type / URI markers
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-comarca"
$ws.Range("F4").Value = "xsd:string"
$ws.Range("G4").Value = "xsd:string"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "null"
